$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.489.63'
$ws.Range('E2').Value = '  +1.43%  '

$ws.Range('D3').Value = '1.668.27'
$ws.Range('E3').Value = '  +1.43%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.46%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.69'
$ws.Range('E5').Value = '  +1.79%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.42%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3970'
$ws.Range('E7').Value = '  +1.62%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3919'
$ws.Range('E8').Value = '  +1.64%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.99'
$ws.Range('E9').Value = '  +6.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.408'
$ws.Range('E10').Value = '  +3.59%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9991'
$ws.Range('E11').Value = '  -0.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08606'
$ws.Range('E12').Value = '  +1.71%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.50'
$ws.Range('E13').Value = '  +1.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.354'
$ws.Range('E14').Value = '  +2.67%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001356'
$ws.Range('E15').Value = '  +5.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.923'
$ws.Range('E16').Value = '  +5.42%  '

$ws.Range('D17').Value = '1.664.08'
$ws.Range('E17').Value = '  +1.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.59'
$ws.Range('E18').Value = '  +1.31%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06988'
$ws.Range('E19').Value = '  +0.66%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.67'
$ws.Range('E20').Value = '  -1.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.029'
$ws.Range('E21').Value = '  +1.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9990'
$ws.Range('E22').Value = '  -0.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.79'
$ws.Range('E23').Value = '  +0.33%  '

$ws.Range('D24').Value = '24.477.76'
$ws.Range('E24').Value = '  +1.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.427'
$ws.Range('E25').Value = '  +3.51%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.047'
$ws.Range('E26').Value = '  +11.47%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.56'
$ws.Range('E27').Value = '  +0.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.03'
$ws.Range('E28').Value = '  -0.02%  '

$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.487'
$ws.Range('E29').Value = '  +1.54%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '143.12'
$ws.Range('E30').Value = '  +1.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.236'
$ws.Range('E31').Value = '  -8.82%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.541'
$ws.Range('E32').Value = '  +3.55%  '

$ws.Range('D33').Value = '1.848.93'
$ws.Range('E33').Value = '  +1.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.071'
$ws.Range('E34').Value = '  +8.78%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08312'
$ws.Range('E35').Value = '  +3.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.03046'
$ws.Range('E36').Value = '  +3.40%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.919'
$ws.Range('E37').Value = '  -3.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.19'
$ws.Range('E38').Value = '  +11.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2786'
$ws.Range('E39').Value = '  +2.82%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09260'
$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.96'
$ws.Range('E41').Value = '  +6.30%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7774'
$ws.Range('E42').Value = '  +1.77%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.449'
$ws.Range('E43').Value = '  -2.02%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.69'
$ws.Range('E44').Value = '  +3.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7155'
$ws.Range('E45').Value = '  +3.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.551'
$ws.Range('E46').Value = '  +2.44%  '

$ws.Range('E47').Value = '  +1.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08469'
$ws.Range('E49').Value = '  +0.60%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.94'
$ws.Range('E50').Value = '  +2.06%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.277'
$ws.Range('E51').Value = '  +0.98%  '

